# Refined metadata to be additional tab
#
# 1) Update the "data" sheet's time_taken (column F) values to the
#    latest query run timestamps.
# 2) Add a new "metadata" worksheet (placed after "data") summarising the
#    PanelApp query that produced this export.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("data")

$timeTakenValues = @(
    "2021-10-05 14:21:04.711076",
    "2021-10-05 14:21:04.711085",
    "2021-10-05 14:21:04.711088",
    "2021-10-05 14:21:04.711091",
    "2021-10-05 14:21:04.711095",
    "2021-10-05 14:21:04.711097",
    "2021-10-05 14:21:04.711100",
    "2021-10-05 14:21:04.711103",
    "2021-10-05 14:21:04.711106",
    "2021-10-05 14:21:04.711109",
    "2021-10-05 14:21:04.711111",
    "2021-10-05 14:21:04.711114",
    "2021-10-05 14:21:04.711117",
    "2021-10-05 14:21:04.711119",
    "2021-10-05 14:21:04.711122",
    "2021-10-05 14:21:04.711124",
    "2021-10-05 14:21:04.711127",
    "2021-10-05 14:21:04.711130",
    "2021-10-05 14:21:04.711133",
    "2021-10-05 14:21:04.711135",
    "2021-10-05 14:21:04.711138",
    "2021-10-05 14:21:04.711141",
    "2021-10-05 14:21:04.711143",
    "2021-10-05 14:21:04.711146",
    "2021-10-05 14:21:04.711149",
    "2021-10-05 14:21:04.711151",
    "2021-10-05 14:21:04.711154",
    "2021-10-05 14:21:04.711157",
    "2021-10-05 14:21:04.711159",
    "2021-10-05 14:21:04.711162",
    "2021-10-05 14:21:04.711165",
    "2021-10-05 14:21:04.711167",
    "2021-10-05 14:21:04.711171",
    "2021-10-05 14:21:04.711173",
    "2021-10-05 14:21:04.711176",
    "2021-10-05 14:21:04.711179",
    "2021-10-05 14:21:04.711181",
    "2021-10-05 14:21:04.711184",
    "2021-10-05 14:21:04.711187",
    "2021-10-05 14:21:04.711189",
    "2021-10-05 14:21:04.711192",
    "2021-10-05 14:21:04.711195",
    "2021-10-05 14:21:04.711198",
    "2021-10-05 14:21:04.711200",
    "2021-10-05 14:21:04.711203",
    "2021-10-05 14:21:04.711205",
    "2021-10-05 14:21:04.711208",
    "2021-10-05 14:21:04.711210",
    "2021-10-05 14:21:04.711213",
    "2021-10-05 14:21:04.711216",
    "2021-10-05 14:21:04.711218",
    "2021-10-05 14:21:04.711221",
    "2021-10-05 14:21:04.711224",
    "2021-10-05 14:21:04.711227",
    "2021-10-05 14:21:04.711229",
    "2021-10-05 14:21:04.711232",
    "2021-10-05 14:21:04.711235",
    "2021-10-05 14:21:04.711237",
    "2021-10-05 14:21:04.711240",
    "2021-10-05 14:21:04.711243",
    "2021-10-05 14:21:04.711246",
    "2021-10-05 14:21:04.711249",
    "2021-10-05 14:21:04.711251",
    "2021-10-05 14:21:04.711254",
    "2021-10-05 14:21:04.711259",
    "2021-10-05 14:21:04.711262",
    "2021-10-05 14:21:04.711264",
    "2021-10-05 14:21:04.711267"
)

for ($i = 0; $i -lt $timeTakenValues.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTakenValues[$i]
}

# --- Add the "metadata" sheet, positioned after "data" ---
$meta = $wb.Worksheets.Add($null, $ws)
$meta.Name = "metadata"

# Copy the header cell style (bold, bordered, centered) from the "data"
# sheet's header row onto the new header row (B1:G1) and the A2 index
# cell, so the new sheet reuses the existing style instead of minting a
# new one.
$ws.Range("B1").Copy()
$meta.Range("B1:G1").PasteSpecial(-4122)
$ws.Range("B1").Copy()
$meta.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$meta.Range("B1").Value = "data_name"
$meta.Range("C1").Value = "data_id"
$meta.Range("D1").Value = "data_version"
$meta.Range("E1").Value = "data_version_created"
$meta.Range("F1").Value = "panel_query_time"
$meta.Range("G1").Value = "panel_get_request"

$meta.Range("A2").Value = 0
$meta.Range("B2").Value = "Infantile enterocolitis & monogenic inflammatory bowel disease"
$meta.Range("C2").Value = 176

# Force the version string to stay text ("1.22") instead of being
# coerced to a number, then drop the number-format override so no new
# style is left behind on the cell.
$meta.Range("D2").NumberFormat = "@"
$meta.Range("D2").Value = "1.22"
$meta.Range("D2").Style = "Normal"

$meta.Range("E2").Value = "2021-08-03T08:06:25.754651Z"
$meta.Range("F2").Value = "2021-10-05 14:21:04.707747"
$meta.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/176/?format=json"

$meta.Range("A1").Select() | Out-Null
$ws.Select() | Out-Null

Write-Output "metadata sheet added; time_taken column refreshed"
